# Applies the Anima_Profits market-data refresh captured in the commit diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets a batch of
# H:N column updates (currentAveragePrice* / LevePrice* / LeveProfit* fields).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2229.8948
$ws.Range("I43").Value = 2900
$ws.Range("J43").Value = 2151.0588
$ws.Range("K43").Value = 2900
$ws.Range("L43").Value = 2151.0588
$ws.Range("M43").Value = -2831
$ws.Range("N43").Value = -2289.0588
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 3182.1177
$ws.Range("I76").Value = 3240
$ws.Range("J76").Value = 3099.4285
$ws.Range("K76").Value = 3240
$ws.Range("L76").Value = 3099.4285
$ws.Range("M76").Value = -2925
$ws.Range("N76").Value = -3729.4285
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 3182.1177
$ws.Range("I79").Value = 3240
$ws.Range("J79").Value = 3099.4285
$ws.Range("K79").Value = 3240
$ws.Range("L79").Value = 3099.4285
$ws.Range("M79").Value = -2148
$ws.Range("N79").Value = -5283.4285
$ws.Range("H140").Value = 77079.336
$ws.Range("J140").Value = 77079.336
$ws.Range("L140").Value = 77079.336
$ws.Range("N140").Value = -87439.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 94677.5
$ws.Range("J24").Value = 94677.5
$ws.Range("L24").Value = 94677.5
$ws.Range("N24").Value = -95425.5
$ws.Range("H32").Value = 22246262
$ws.Range("I32").Value = 35737028
$ws.Range("K32").Value = 35737028
$ws.Range("M32").Value = -35736741
$ws.Range("H37").Value = 22464.666
$ws.Range("I37").Value = 2034
$ws.Range("K37").Value = 2034
$ws.Range("M37").Value = -1761
$ws.Range("H45").Value = 2097
$ws.Range("I45").Value = 999
$ws.Range("K45").Value = 999
$ws.Range("M45").Value = -622
$ws.Range("H100").Value = 94677.5
$ws.Range("J100").Value = 94677.5
$ws.Range("L100").Value = 94677.5
$ws.Range("N100").Value = -96841.5
$ws.Range("H122").Value = 201222.8
$ws.Range("I122").Value = 251150
$ws.Range("K122").Value = 753450
$ws.Range("M122").Value = -751000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 475.8
$ws.Range("I64").Value = 512.625
$ws.Range("J64").Value = 433.7143
$ws.Range("K64").Value = 512.625
$ws.Range("L64").Value = 433.7143
$ws.Range("M64").Value = -287.625
$ws.Range("N64").Value = -883.7143
$ws.Range("H67").Value = 475.8
$ws.Range("I67").Value = 512.625
$ws.Range("J67").Value = 433.7143
$ws.Range("K67").Value = 512.625
$ws.Range("L67").Value = 433.7143
$ws.Range("M67").Value = 267.375
$ws.Range("N67").Value = -1993.7143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 47004.668
$ws.Range("I2").Value = 1004
$ws.Range("J2").Value = 70005
$ws.Range("K2").Value = 1004
$ws.Range("L2").Value = 70005
$ws.Range("M2").Value = -891
$ws.Range("N2").Value = -70231
$ws.Range("I7").Value = 101
$ws.Range("J7").Value = 64.40000000000001
$ws.Range("K7").Value = 101
$ws.Range("L7").Value = 64.40000000000001
$ws.Range("M7").Value = 12
$ws.Range("N7").Value = -290.4
$ws.Range("H31").Value = 7471.1816
$ws.Range("I31").Value = 3088
$ws.Range("J31").Value = 8651.27
$ws.Range("K31").Value = 3088
$ws.Range("L31").Value = 8651.27
$ws.Range("M31").Value = -2793
$ws.Range("N31").Value = -9241.27
$ws.Range("H34").Value = 7471.1816
$ws.Range("I34").Value = 3088
$ws.Range("J34").Value = 8651.27
$ws.Range("K34").Value = 3088
$ws.Range("L34").Value = 8651.27
$ws.Range("M34").Value = -2886
$ws.Range("N34").Value = -9055.27

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7369346.5
$ws.Range("I4").Value = 10000470
$ws.Range("J4").Value = 2199.8
$ws.Range("K4").Value = 30001410
$ws.Range("L4").Value = 6599.400000000001
$ws.Range("M4").Value = -30001298
$ws.Range("N4").Value = -6823.400000000001
$ws.Range("H131").Value = 2811.2031
$ws.Range("I131").Value = 608.3077
$ws.Range("J131").Value = 3372.7256
$ws.Range("K131").Value = 1824.9231
$ws.Range("L131").Value = 10118.1768
$ws.Range("M131").Value = 3215.0769
$ws.Range("N131").Value = -20198.1768
$ws.Range("H137").Value = 9840817
$ws.Range("I137").Value = 11039.8
$ws.Range("J137").Value = 17601168
$ws.Range("K137").Value = 33119.39999999999
$ws.Range("L137").Value = 52803504
$ws.Range("M137").Value = -28019.39999999999
$ws.Range("N137").Value = -52813704
$ws.Range("H140").Value = 1732.95
$ws.Range("I140").Value = 1465.1052
$ws.Range("J140").Value = 1975.2858
$ws.Range("K140").Value = 4395.3156
$ws.Range("L140").Value = 5925.857400000001
$ws.Range("M140").Value = 784.6844000000001
$ws.Range("N140").Value = -16285.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4744.3823
$ws.Range("I70").Value = 4718.2085
$ws.Range("K70").Value = 4718.2085
$ws.Range("M70").Value = -4448.2085
$ws.Range("H73").Value = 4744.3823
$ws.Range("I73").Value = 4718.2085
$ws.Range("K73").Value = 4718.2085
$ws.Range("M73").Value = -3782.2085
$ws.Range("H122").Value = 1736.25
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 1555.7142
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 4667.142599999999
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -9567.142599999999
$ws.Range("H123").Value = 16463.234
$ws.Range("J123").Value = 16463.234
$ws.Range("L123").Value = 16463.234
$ws.Range("N123").Value = -21363.234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9496
$ws.Range("J2").Value = 9496
$ws.Range("L2").Value = 9496
$ws.Range("N2").Value = -9720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1000000000
$ws.Range("I4").Value = 1000000000
$ws.Range("K4").Value = 1000000000
$ws.Range("M4").Value = -999999887
